$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row rename (A1:D1) ---
$ws.Range("A1").Value = "mx_state"
$ws.Range("B1").Value = "mx_municipality"
$ws.Range("C1").Value = "n_matriculas"
$ws.Range("D1").Value = "pct_matriculas"

# --- Column A (state name) "de"/"del" -> "De"/"Del" title-casing ---
$colA = @{
    197 = 'Ciudad De México'
    248 = 'Estado De México'
}
foreach ($r in $colA.Keys) {
    $ws.Cells.Item($r, 1).Value = $colA[$r]
}

# --- Column B (municipality name) "de"/"del"/"la"/"los"/"el"/"y" -> title-cased connector words ---
$colB = @{
    8 = 'Pabellón De Arteaga'
    9 = 'Rincón De Romos'
    10 = 'San Francisco De Los Romo'
    11 = 'San José De Gracia'
    39 = 'Amatenango De La Frontera'
    40 = 'Amatenango Del Valle'
    43 = 'Bejucal De Ocampo'
    53 = 'Chiapa De Corzo'
    58 = 'Comitán De Domínguez'
    85 = 'Marqués De Comillas'
    86 = 'Mazapa De Madero'
    92 = 'Ocozocoautla De Espinosa'
    104 = 'Salto De Agua'
    105 = 'San Cristóbal De Las Casas'
    143 = 'Guadalupe Y Calvo'
    144 = 'Hidalgo Del Parral'
    154 = 'San Francisco Del Oro'
    158 = 'Valle De Zaragoza'
    181 = 'San Juan De Sabinas'
    195 = 'Villa De Álvarez'
    201 = 'Cuajimalpa De Morelos'
    215 = 'Coneto De Comonfort'
    229 = 'Nombre De Dios'
    233 = 'Pánuco De Coronado'
    239 = 'San Juan De Guadalupe'
    240 = 'San Juan Del Río'
    248 = 'Acambay De Ruíz Castañeda'
    251 = 'Almoloya De Alquisiras'
    252 = 'Almoloya De Juárez'
    253 = 'Almoloya Del Río'
    259 = 'Atizapán De Zaragoza'
    264 = 'Chapa De Mota'
    268 = 'Coacalco De Berriozábal'
    274 = 'Ecatepec De Morelos'
    281 = 'Ixtapan De La Sal'
    282 = 'Ixtapan Del Oro'
    297 = 'Naucalpan De Juárez'
    309 = 'San Felipe Del Progreso'
    311 = 'San Simón De Guerrero'
    313 = 'Soyaniquilpan De Juárez'
    321 = 'Tenango Del Valle'
    332 = 'Tlalnepantla De Baz'
    338 = 'Valle De Bravo'
    339 = 'Valle De Chalco Solidaridad'
    340 = 'Villa De Allende'
    341 = 'Villa Del Carbón'
    353 = 'San Miguel De Allende'
    354 = 'Apaseo El Alto'
    355 = 'Apaseo El Grande'
    363 = 'Dolores Hidalgo Cuna De La Independencia Nacional'
    367 = 'Jaral Del Progreso'
    375 = 'Purísima Del Rincón'
    379 = 'San Diego De La Unión'
    381 = 'San Francisco Del Rincón'
    383 = 'San Luis De La Paz'
    385 = 'Santa Cruz De Juventino Rosas'
    387 = 'Silao De La Victoria'
    392 = 'Valle De Santiago'
    398 = 'Acapulco De Juárez'
    401 = 'Ajuchitlán Del Progreso'
    402 = 'Alcozauca De Guerrero'
    406 = 'Atenango Del Río'
    407 = 'Atlamajalcingo Del Monte'
    409 = 'Atoyac De Álvarez'
    410 = 'Ayutla De Los Libres'
    413 = 'Buenavista De Cuéllar'
    414 = 'Chilapa De Álvarez'
    415 = 'Chilpancingo De Los Bravo'
    416 = 'Coahuayutla De José María Izazaga'
    421 = 'Coyuca De Benítez'
    422 = 'Coyuca De Catalán'
    426 = 'Cuetzala Del Progreso'
    427 = 'Cutzamala De Pinzón'
    433 = 'Huitzuco De Los Figueroa'
    434 = 'Iguala De La Independencia'
    436 = 'Ixcateopan De Cuauhtémoc'
    437 = 'Zihuatanejo De Azueta'
    439 = 'La Unión De Isidoro Montes De Oca'
    442 = 'Mártir De Cuilapan'
    455 = 'Taxco De Alarcón'
    457 = 'Técpan De Galeana'
    459 = 'Tepecoacuilco De Trujano'
    461 = 'Tixtla De Guerrero'
    465 = 'Tlalixtaquilla De Maldonado'
    466 = 'Tlapa De Comonfort'
    478 = 'Agua Blanca De Iturbide'
    485 = 'Atotonilco De Tula'
    486 = 'Atotonilco El Grande'
    492 = 'Cuautepec De Hinojosa'
    496 = 'Huasca De Ocampo'
    500 = 'Huejutla De Reyes'
    503 = 'Jacala De Ledezma'
    510 = 'Mineral De La Reforma'
    511 = 'Mineral Del Chico'
    512 = 'Mineral Del Monte'
    513 = 'Mixquiahuala De Juárez'
    514 = 'Molango De Escamilla'
    516 = 'Nopala De Villagrán'
    517 = 'Omitlán De Juárez'
    518 = 'Pachuca De Soto'
    521 = 'Progreso De Obregón'
    527 = 'Santiago De Anaya'
    528 = 'Santiago Tulantepec De Lugo Guerrero'
    532 = 'Tenango De Doria'
    534 = 'Tepehuacán De Guerrero'
    535 = 'Tepeji Del Río De Ocampo'
    538 = 'Tezontepec De Aldama'
    545 = 'Tula De Allende'
    546 = 'Tulancingo De Bravo'
    547 = 'Villa De Tezontepec'
    550 = 'Zacualtipán De Ángeles'
    555 = 'Acatlán De Juárez'
    556 = 'Ahualulco De Mercado'
    561 = 'Atotonilco El Alto'
    563 = 'Autlán De Navarro'
    568 = 'Cañadas De Obregón'
    579 = 'Encarnación De Díaz'
    583 = 'Ixtlahuacán Del Río'
    587 = 'Jilotlán De Los Dolores'
    592 = 'Lagos De Moreno'
    598 = 'Ojuelos De Jalisco'
    603 = 'San Diego De Alejandría'
    605 = 'San Juan De Los Lagos'
    608 = 'San Miguel El Alto'
    609 = 'San Sebastián Del Oeste'
    610 = 'Santa María De Los Ángeles'
    611 = 'Santa María Del Oro'
    614 = 'Talpa De Allende'
    615 = 'Tamazula De Gordiano'
    617 = 'Techaluta De Montenegro'
    621 = 'Teocuitatlán De Corona'
    622 = 'Tepatitlán De Morelos'
    624 = 'Tizapán El Alto'
    633 = 'Unión De San Antonio'
    634 = 'Unión De Tula'
    635 = 'Valle De Juárez'
    639 = 'Yahualica De González Gallo'
    640 = 'Zacoalco De Torres'
    643 = 'Zapotlán Del Rey'
    644 = 'Zapotlán El Grande'
    668 = 'Coalcomán De Vázquez Pallares'
    733 = 'Tiquicheo De Nicolás Romero'
    758 = 'Coatlán Del Río'
    765 = 'Jonacatepec De Leandro Valle'
    769 = 'Puente De Ixtla'
    773 = 'Tetela Del Volcán'
    775 = 'Tlaltizapán De Zapata'
    782 = 'Zacualpan De Amilpas'
    785 = 'Bahía De Banderas'
    789 = 'Ixtlán Del Río'
    796 = 'Santa María Del Oro'
    815 = 'Lampazos De Naranjo'
    817 = 'Mier Y Noriega'
    822 = 'San Nicolás De Los Garza'
    826 = 'Acatlán De Pérez Figueroa'
    833 = 'Ayoquezco De Aldama'
    836 = 'Chalcatongo De Hidalgo'
    837 = 'Ciénega De Zimatlán'
    839 = 'Coicoyán De Las Flores'
    840 = 'Constancia Del Rosario'
    843 = 'Eloxochitlán De Flores Magón'
    844 = 'Fresnillo De Trujano'
    845 = 'Guadalupe De Ramírez'
    847 = 'Guelatao De Juárez'
    848 = 'Heroica Ciudad De Ejutla De Crespo'
    849 = 'Heroica Ciudad De Huajuapan De León'
    850 = 'Heroica Ciudad De Tlaxiaco'
    851 = 'Huautla De Jiménez'
    853 = 'Ixtlán De Juárez'
    854 = 'Heroica Ciudad De Juchitán De Zaragoza'
    862 = 'Mariscala De Juárez'
    863 = 'Mártires De Tacubaya'
    865 = 'Mazatlán Villa De Flores'
    867 = 'Miahuatlán De Porfirio Díaz'
    871 = 'Oaxaca De Juárez'
    872 = 'Ocotlán De Morelos'
    873 = 'Pinotepa De Don Luis'
    875 = 'Putla Villa De Guerrero'
    876 = 'Reforma De Pineda'
    889 = 'San Antonino El Alto'
    891 = 'San Antonio De La Cal'
    898 = 'San Baltazar Yatzachi El Bajo'
    909 = 'San Felipe Jalapa De Díaz'
    913 = 'San Francisco Del Mar'
    929 = 'San José Del Progreso'
    944 = 'San Juan Del Estado'
    945 = 'San Juan Del Río'
    986 = 'San Miguel Del Puerto'
    988 = 'San Miguel El Grande'
    1004 = 'San Pablo Villa De Mitla'
    1009 = 'San Pedro El Alto'
    1025 = 'San Pedro Y San Pablo Ayutla'
    1026 = 'San Pedro Y San Pablo Teposcolula'
    1027 = 'San Pedro Y San Pablo Tequixtepec'
    1046 = 'Santa Cruz De Bravo'
    1050 = 'Santa Cruz Tacache De Mina'
    1055 = 'Santa Inés Del Monte'
    1057 = 'Santa Lucía Del Camino'
    1073 = 'Santa María Jalapa Del Marqués'
    1095 = 'Santiago Del Río'
    1122 = 'Santo Domingo De Morelos'
    1137 = 'Tamazulápam Del Espíritu Santo'
    1139 = 'Tataltepec De Valdés'
    1140 = 'Teotitlán De Flores Magón'
    1141 = 'Teotitlán Del Valle'
    1143 = 'Tepelmeme Villa De Morelos'
    1144 = 'Heroica Villa Tezoatlán De Segura Y Luna, Cuna De La Independencia De Oaxaca'
    1145 = 'Tlacolula De Matamoros'
    1146 = 'Totontepec Villa De Morelos'
    1148 = 'Villa De Chilapa De Díaz'
    1149 = 'Villa De Etla'
    1150 = 'Villa De Tamazulápam Del Progreso'
    1151 = 'Villa De Tututepec'
    1152 = 'Villa De Zaachila'
    1155 = 'Villa Sola De Vega'
    1156 = 'Villa Talea De Castro'
    1157 = 'Villa Tejúpam De La Unión'
    1160 = 'Zapotitlán Del Río'
    1163 = 'Zimatlán De Álvarez'
    1182 = 'Ayotoxco De Guerrero'
    1184 = 'Chalchicomula De Sesma'
    1193 = 'Chila De La Sal'
    1203 = 'Cuetzalan Del Progreso'
    1215 = 'Huehuetlán El Chico'
    1219 = 'Huitzilan De Serdán'
    1220 = 'Ixcamilpa De Guerrero'
    1223 = 'Izúcar De Matamoros'
    1233 = 'Los Reyes De Juárez'
    1234 = 'Mazapiltepec De Juárez'
    1244 = 'Palmar De Bravo'
    1261 = 'San Nicolás De Los Ranchos'
    1264 = 'San Salvador El Seco'
    1265 = 'San Salvador El Verde'
    1275 = 'Tepanco De López'
    1276 = 'Tepatlaxco De Hidalgo'
    1280 = 'Tepexi De Rodríguez'
    1282 = 'Tetela De Ocampo'
    1287 = 'Tlacotepec De Benito Juárez'
    1298 = 'Tuzamapan De Galeana'
    1302 = 'Xayacatlán De Bravo'
    1312 = 'Zapotitlán De Méndez'
    1319 = 'Amealco De Bonfil'
    1321 = 'Cadereyta De Montes'
    1327 = 'Jalpan De Serra'
    1328 = 'Landa De Matamoros'
    1331 = 'Pinal De Amoles'
    1334 = 'San Juan Del Río'
    1345 = 'Armadillo De Los Infante'
    1346 = 'Axtla De Terrazas'
    1351 = 'Cerro De San Pedro'
    1353 = 'Ciudad Del Maíz'
    1362 = 'Mexquitic De Carmona'
    1367 = 'San Ciro De Acosta'
    1372 = 'Santa María Del Río'
    1378 = 'Tanquián De Escobedo'
    1382 = 'Villa De Arista'
    1383 = 'Villa De Arriaga'
    1384 = 'Villa De Guadalupe'
    1385 = 'Villa De La Paz'
    1386 = 'Villa De Ramos'
    1387 = 'Villa De Reyes'
    1432 = 'Jalpa De Méndez'
    1472 = 'Soto La Marina'
    1485 = 'Contla De Juan Cuamatzi'
    1490 = 'Ixtacuixtla De Mariano Matamoros'
    1492 = 'San Pablo Del Monte'
    1513 = 'Alto Lucero De Gutiérrez Barrios'
    1517 = 'Amatlán De Los Reyes'
    1528 = 'Boca Del Río'
    1530 = 'Camarón De Tejeda'
    1534 = 'Castillo De Teayo'
    1536 = 'Cazones De Herrera'
    1544 = 'Chinampa De Gorostiza'
    1553 = 'Cosamaloapan De Carpio'
    1569 = 'Hueyapan De Ocampo'
    1570 = 'Ignacio De La Llave'
    1573 = 'Ixhuatlán De Madero'
    1574 = 'Ixhuatlán Del Café'
    1575 = 'Ixhuatlán Del Sureste'
    1585 = 'Juchique De Ferrer'
    1588 = 'Landero Y Coss'
    1590 = 'Las Vigas De Ramírez'
    1591 = 'Lerdo De Tejada'
    1595 = 'Martínez De La Torre'
    1597 = 'Medellín De Bravo'
    1601 = 'Mixtla De Altamirano'
    1603 = 'Nanchital De Lázaro Cárdenas Del Río'
    1613 = 'Ozuluama De Mascareñas'
    1617 = 'Paso De Ovejas'
    1618 = 'Paso Del Macho'
    1622 = 'Poza Rica De Hidalgo'
    1631 = 'Sayula De Alemán'
    1633 = 'Soledad De Doblado'
    1640 = 'Tatahuicapan De Juárez'
    1670 = 'Vega De Alatorre'
    1679 = 'Zontecomatlán De López Y Fuentes'
    1696 = 'Cañitas De Felipe Pescador'
    1698 = 'Concepción Del Oro'
    1707 = 'Jiménez Del Teul'
    1716 = 'Moyahua De Estrada'
    1717 = 'Nochistlán De Mejía'
    1718 = 'Noria De Ángeles'
    1726 = 'Tlaltenango De Sánchez Román'
    1730 = 'Villa De Cos'
}
foreach ($r in $colB.Keys) {
    $ws.Cells.Item($r, 2).Value = $colB[$r]
}

# --- Floating point precision fixes (1 ULP) on column D ---
$colD = @{
    78 = 0.0009721952167995332
    467 = 0.0009721952167995332
    492 = 0.0009721952167995332
    1163 = 0.0009721952167995332
    1164 = 0.09313630176939527
    1361 = 0.0009721952167995332
    1600 = 0.0009721952167995332
}
foreach ($r in $colD.Keys) {
    $ws.Cells.Item($r, 4).Value = $colD[$r]
}

# --- Remove trailing footnote rows 1739-1743 (dimension shrinks to A1:D1737) ---
$ws.Rows("1739:1743").Delete()

